$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.239.41"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.49%  "

# Row 3: Ethereum
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.990.15"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +5.99%  "

# Row 4: TetherUSD (price unchanged)
$ws.Range("E4").Value = "  -0.26%  "

# Row 5: BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "325.35"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.77%  "

# Row 6: USDC (price unchanged)
$ws.Range("E6").Value = "  -0.21%  "

# Row 7: XRP
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5103"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.32%  "

# Row 8: Cardano
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4149"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.78%  "

# Row 9: Dogecoin
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08703"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.97%  "

# Row 10: Polygon
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.130"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.38%  "

# Row 11: OKB
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.67"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.43%  "

# Row 12: Solana
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.21"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.71%  "

# Row 13: WrappedEther
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.994.60"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +6.80%  "

# Row 14: Polkadot
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.492"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.12%  "

# Row 15: Chainlink
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.389"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.71%  "

# Row 16: BinanceUSD
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.002"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.24%  "

# Row 17: Litecoin
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "93.93"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.40%  "

# Row 18: ShibaInu (price unchanged)
$ws.Range("E18").Value = "  +2.43%  "

# Row 19: TRON
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06552"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.47%  "

# Row 20: Avalanche
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.77"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.83%  "

# Row 21: Dai (price unchanged)
$ws.Range("E21").Value = "  -0.14%  "

# Row 22: Uniswap
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.068"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.94%  "

# Row 23: WrappedBTC
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.310.65"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.74%  "

# Row 24: Cosmos
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.58"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.80%  "

# Row 25: Toncoin
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.203"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.30%  "

# Row 26: WrappedliquidstakedEther2.0
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.224.59"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +6.61%  "

# Row 27: EthereumClassic
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.52"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +6.64%  "

# Row 28: Monero
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "163.10"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.17%  "

# Row 29: LidoDAOToken
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.372"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.56%  "

# Row 30: BitcoinCash
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "130.70"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.45%  "

# Row 31: ImmutableX
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.130"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.25%  "

# Row 32: Stellar (price unchanged)
$ws.Range("E32").Value = "  +1.62%  "

# Row 33: Filecoin
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.040"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.00%  "

# Row 34: HuobiToken
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.812"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.12%  "

# Row 35: ARBITRUM (price unchanged)
$ws.Range("E35").Value = "  +12.02%  "

# Row 36: VeChain
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02479"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.20%  "

# Row 37: InternetComputer(DFINITY)
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.385"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.28%  "

# Row 38: Hedera
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06516"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.65%  "

# Row 39: Algorand
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2192"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.95%  "

# Row 40: FraxShare
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.895"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.76%  "

# Row 41: TheSandbox
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6567"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.43%  "

# Row 42: Aptos
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.76"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.26%  "

# Row 43: TrustWalletToken (price unchanged)
$ws.Range("E43").Value = "  +1.02%  "

# Row 44: EnergySwap
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.67"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.91%  "

# Row 45: Decentraland
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6116"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.59%  "

# Row 46: NEARProtocol
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.197"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.24%  "

# Row 47: PancakeSwap
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.666"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.05%  "

# Row 48: Quant
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "124.12"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.71%  "

# Row 49: EOS
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.224"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.27%  "

# Row 50: Aave
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "79.34"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.57%  "

# Row 51: Cronos
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06871"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.01%  "
